$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.28
$ws.Range("I2").Value = 3.45
$ws.Range("K2").Value = 4
$ws.Range("V2").Value = 1.41
$ws.Range("W2").Value = 1.78

# Row 4
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.75

# Row 5
$ws.Range("I5").Value = 2.2
$ws.Range("J5").Value = 3.5
$ws.Range("P5").Value = 1.95
$ws.Range("Q5").Value = 1.84
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 2.88
$ws.Range("V5").Value = 1.89
$ws.Range("W5").Value = 1.24

# Row 7
$ws.Range("G7").Value = 840
$ws.Range("I7").Value = 730
$ws.Range("K7").Value = 95

# Row 8
$ws.Range("G8").Value = 1.4
$ws.Range("L8").Value = 1.01
$ws.Range("V8").Value = 1.1
$ws.Range("W8").Value = 3.5
$ws.Range("AB8").Value = 1000
$ws.Range("AE8").Value = 130
$ws.Range("AI8").Value = 90

# Row 12
$ws.Range("G12").Value = 2.06
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 4.1
$ws.Range("J12").Value = 4.4
$ws.Range("K12").Value = 5.1
$ws.Range("Q12").Value = 1.37

# Row 13
$ws.Range("H13").Value = 2.6
$ws.Range("P13").Value = 2.8
$ws.Range("Q13").Value = 1.45
$ws.Range("R13").Value = 1.68
$ws.Range("T13").Value = 1.41

# Row 14
$ws.Range("I14").Value = 4.8
$ws.Range("K14").Value = 5.1

# Row 16
$ws.Range("F16").Value = 1.41
$ws.Range("G16").Value = 1.51
$ws.Range("H16").Value = 6.6
$ws.Range("I16").Value = 8.800000000000001
$ws.Range("J16").Value = 5.1
$ws.Range("K16").Value = 6.8
$ws.Range("P16").Value = 2.9

# Row 17
$ws.Range("G17").Value = 3.05
$ws.Range("I17").Value = 5.1
$ws.Range("J17").Value = 2.82
$ws.Range("P17").Value = 1.53

# Row 20
$ws.Range("F20").Value = 1.95
$ws.Range("G20").Value = 2.2
$ws.Range("H20").Value = 4
$ws.Range("I20").Value = 4.7
$ws.Range("J20").Value = 3.35
$ws.Range("K20").Value = 3.8

# Row 21
$ws.Range("F21").Value = 2.58
$ws.Range("G21").Value = 2.9
$ws.Range("I21").Value = 3.35

# Row 22
$ws.Range("F22").Value = 3
$ws.Range("I22").Value = 2.92
$ws.Range("L22").Value = 1.47
$ws.Range("V22").Value = 1.52
$ws.Range("W22").Value = 1.48

# Row 23
$ws.Range("F23").Value = 1.27
$ws.Range("G23").Value = 3.05
$ws.Range("H23").Value = 3.45
$ws.Range("I23").Value = 4.6
$ws.Range("J23").Value = 2.18

# Row 24
$ws.Range("G24").Value = 2.64
$ws.Range("H24").Value = 3.1
$ws.Range("L24").Value = 1.01
$ws.Range("R24").Value = 1.28
$ws.Range("V24").Value = 1.45
$ws.Range("W24").Value = 1.61

# Row 25
$ws.Range("M25").Value = 1.15
$ws.Range("AE25").Value = 1000

# Row 26
$ws.Range("F26").Value = 1.38
$ws.Range("I26").Value = 13
$ws.Range("K26").Value = 5.3
$ws.Range("P26").Value = 1.81

# Row 27
$ws.Range("Q27").Value = 2.1
